$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 20.56721760780198
$ws.Range("C2").Value = 3.450500540522718
$ws.Range("D2").Value = 4.067043080912397
$ws.Range("F2").Value = 70.59954383540085
$ws.Range("G2").Value = 3.844794244299919
$ws.Range("J2").Value = 11.52536967609948
$ws.Range("K2").Value = 16.77824811752698
$ws.Range("L2").Value = 11.0516073186643
$ws.Range("M2").Value = 18.64989958963955
$ws.Range("B3").Value = 20.59140359385286
$ws.Range("C3").Value = 3.31192558407896
$ws.Range("D3").Value = 3.930455569181884
$ws.Range("F3").Value = 69.77525348130104
$ws.Range("G3").Value = 3.848371187880647
$ws.Range("J3").Value = 11.51528772533953
$ws.Range("K3").Value = 16.81308973338104
$ws.Range("L3").Value = 11.09308149487794
$ws.Range("M3").Value = 18.71376382645267
$ws.Range("B4").Value = 20.61255952927002
$ws.Range("C4").Value = 3.222812220202493
$ws.Range("D4").Value = 3.843166324378965
$ws.Range("F4").Value = 69.26632504651096
$ws.Range("G4").Value = 3.850680100466366
$ws.Range("J4").Value = 11.50927960021175
$ws.Range("K4").Value = 16.83953969445869
$ws.Range("L4").Value = 11.12017425581475
$ws.Range("M4").Value = 18.75667504861791
$ws.Range("B5").Value = 20.62276370249883
$ws.Range("C5").Value = 3.185501894069224
$ws.Range("D5").Value = 3.806757267809792
$ws.Range("F5").Value = 69.05836091342678
$ws.Range("G5").Value = 3.851649437185744
$ws.Range("J5").Value = 11.50687711971054
$ws.Range("K5").Value = 16.85158673916519
$ws.Range("L5").Value = 11.13162492802477
$ws.Range("M5").Value = 18.77509137806127
$ws.Range("B6").Value = 20.62455361063181
$ws.Range("C6").Value = 3.179246938819651
$ws.Range("D6").Value = 3.800661682892788
$ws.Range("F6").Value = 69.02379804218936
$ws.Range("G6").Value = 3.851812115454563
$ws.Range("J6").Value = 11.50648096206295
$ws.Range("K6").Value = 16.85366365224009
$ws.Range("L6").Value = 11.13355110269251
$ws.Range("M6").Value = 18.77820553496748
$ws.Range("B7").Value = 20.61269074117126
$ws.Range("C7").Value = 3.222313048160878
$ws.Range("D7").Value = 3.842678659717425
$ws.Range("F7").Value = 69.26352250609447
$ws.Range("G7").Value = 3.850693058002461
$ws.Range("J7").Value = 11.50924701367717
$ws.Range("K7").Value = 16.83969703373071
$ws.Range("L7").Value = 11.12032702168717
$ws.Range("M7").Value = 18.75691965400103
$ws.Range("B8").Value = 20.57424685581486
$ws.Range("C8").Value = 3.403565201205232
$ws.Range("D8").Value = 4.020671352054241
$ws.Range("F8").Value = 70.31599298491962
$ws.Range("G8").Value = 3.846004258187846
$ws.Range("J8").Value = 11.52185562992348
$ws.Range("K8").Value = 16.78921014706729
$ws.Range("L8").Value = 11.0655703878617
$ws.Range("M8").Value = 18.67115220911004
$ws.Range("B9").Value = 20.54897325132769
$ws.Range("C9").Value = 3.726428147134247
$ws.Range("D9").Value = 4.341666332362025
$ws.Range("F9").Value = 72.35167359064224
$ws.Range("G9").Value = 3.837698464383422
$ws.Range("J9").Value = 11.54802562892953
$ws.Range("K9").Value = 16.73045258309247
$ws.Range("L9").Value = 10.97106476426886
$ws.Range("M9").Value = 18.53232096406559
$ws.Range("B10").Value = 20.56101505536672
$ws.Range("C10").Value = 3.943118649907817
$ws.Range("D10").Value = 4.559392485985184
$ws.Range("F10").Value = 73.82226292740151
$ws.Range("G10").Value = 3.832131197531965
$ws.Range("J10").Value = 11.56813302274881
$ws.Range("K10").Value = 16.71195475266283
$ws.Range("L10").Value = 10.9094220962124
$ws.Range("M10").Value = 18.44823804266672
$ws.Range("B11").Value = 20.57313340693244
$ws.Range("C11").Value = 4.03714105011611
$ws.Range("D11").Value = 4.654339116989576
$ws.Range("F11").Value = 74.48420854127289
$ws.Range("G11").Value = 3.829713196096534
$ws.Range("J11").Value = 11.57747136710586
$ws.Range("K11").Value = 16.70891455508393
$ws.Range("L11").Value = 10.88305909906587
$ws.Range("M11").Value = 18.4138805459841
$ws.Range("B12").Value = 20.57867500078321
$ws.Range("C12").Value = 4.072083970373098
$ws.Range("D12").Value = 4.689691786349909
$ws.Range("F12").Value = 74.73372464045225
$ws.Range("G12").Value = 3.828813924326566
$ws.Range("J12").Value = 11.58103498997018
$ws.Range("K12").Value = 16.70853663084583
$ws.Range("L12").Value = 10.87331660764635
$ws.Range("M12").Value = 18.40143023675516
$ws.Range("B13").Value = 20.57743919696922
$ws.Range("C13").Value = 4.064587894482044
$ws.Range("D13").Value = 4.68210488819432
$ws.Range("F13").Value = 74.68004010386188
$ws.Range("G13").Value = 3.829006872039133
$ws.Range("J13").Value = 11.58026628380106
$ws.Range("K13").Value = 16.70858362823897
$ws.Range("L13").Value = 10.87540413871853
$ws.Range("M13").Value = 18.40408671525482
$ws.Range("B14").Value = 20.57357024134564
$ws.Range("C14").Value = 4.040029113494914
$ws.Range("D14").Value = 4.657259722701243
$ws.Range("F14").Value = 74.50476016869278
$ws.Range("G14").Value = 3.82963888491875
$ws.Range("J14").Value = 11.57776400024586
$ws.Range("K14").Value = 16.70886796370972
$ws.Range("L14").Value = 10.88225276075231
$ws.Range("M14").Value = 18.41284502339959
$ws.Range("B15").Value = 20.57132437205031
$ws.Range("C15").Value = 4.024899847266148
$ws.Range("D15").Value = 4.641962665922815
$ws.Range("F15").Value = 74.39724258792202
$ws.Range("G15").Value = 3.830028140240317
$ws.Range("J15").Value = 11.57623484291878
$ws.Range("K15").Value = 16.70914284197712
$ws.Range("L15").Value = 10.8864790505072
$ws.Range("M15").Value = 18.4182826974986
$ws.Range("B16").Value = 20.56035654562072
$ws.Range("C16").Value = 3.93688186418665
$ws.Range("D16").Value = 4.553103949082576
$ws.Range("F16").Value = 73.77884992085174
$ws.Range("G16").Value = 3.83229151636259
$ws.Range("J16").Value = 11.56752658064651
$ws.Range("K16").Value = 16.71226164118249
$ws.Range("L16").Value = 10.91117869073766
$ws.Range("M16").Value = 18.45056174526771
$ws.Range("B17").Value = 20.5553276319562
$ws.Range("C17").Value = 3.881714674201192
$ws.Range("D17").Value = 4.497533111404743
$ws.Range("F17").Value = 73.39758629478318
$ws.Range("G17").Value = 3.833709296575864
$ws.Range("J17").Value = 11.56223327673663
$ws.Range("K17").Value = 16.71555200767268
$ws.Range("L17").Value = 10.92676048630379
$ws.Range("M17").Value = 18.47136110911618
$ws.Range("B18").Value = 20.55306036423506
$ws.Range("C18").Value = 3.849555899823631
$ws.Range("D18").Value = 4.465185121437658
$ws.Range("F18").Value = 73.17763958983971
$ws.Range("G18").Value = 3.834535556843778
$ws.Range("J18").Value = 11.55920673123873
$ws.Range("K18").Value = 16.71795042784567
$ws.Range("L18").Value = 10.93588076245412
$ws.Range("M18").Value = 18.48369068628653
$ws.Range("B19").Value = 20.55240014278284
$ws.Range("C19").Value = 3.838594203120074
$ws.Range("D19").Value = 4.454166906087702
$ws.Range("F19").Value = 73.10306132408051
$ws.Range("G19").Value = 3.834817170811466
$ws.Range("J19").Value = 11.55818509407107
$ws.Range("K19").Value = 16.71884935403671
$ws.Range("L19").Value = 10.93899589840531
$ws.Range("M19").Value = 18.48792817366405
$ws.Range("B20").Value = 20.55579827280392
$ws.Range("C20").Value = 3.887631673005199
$ws.Range("D20").Value = 4.503488673899182
$ws.Range("F20").Value = 73.43824114291431
$ws.Range("G20").Value = 3.833557255418713
$ws.Range("J20").Value = 11.56279489353336
$ws.Range("K20").Value = 16.71514937978407
$ws.Range("L20").Value = 10.92508542676549
$ws.Range("M20").Value = 18.46910906412908
$ws.Range("B21").Value = 20.5746808159484
$ws.Range("C21").Value = 4.047260618753361
$ws.Range("D21").Value = 4.664573765495493
$ws.Range("F21").Value = 74.55627637602613
$ws.Range("G21").Value = 3.82945280387909
$ws.Range("J21").Value = 11.57849823917585
$ws.Range("K21").Value = 16.70876345865512
$ws.Range("L21").Value = 10.88023463061802
$ws.Range("M21").Value = 18.41025729100388
$ws.Range("B22").Value = 20.59257260653813
$ws.Range("C22").Value = 4.147730426233403
$ws.Range("D22").Value = 4.766341316517289
$ws.Range("F22").Value = 75.2802224563446
$ws.Range("G22").Value = 3.826865693774279
$ws.Range("J22").Value = 11.5889208922814
$ws.Range("K22").Value = 16.70909725675564
$ws.Range("L22").Value = 10.85232411566752
$ws.Range("M22").Value = 18.37505929005197
$ws.Range("B23").Value = 20.58251650748877
$ws.Range("C23").Value = 4.094462628061547
$ws.Range("D23").Value = 4.712350922499212
$ws.Range("F23").Value = 74.89450126275527
$ws.Range("G23").Value = 3.828237789208919
$ws.Range("J23").Value = 11.58334356699546
$ws.Range("K23").Value = 16.70850667695965
$ws.Range("L23").Value = 10.86709244422874
$ws.Range("M23").Value = 18.39354624207727
$ws.Range("B24").Value = 20.55558355245772
$ws.Range("C24").Value = 3.884957975409139
$ws.Range("D24").Value = 4.500797407792927
$ws.Range("F24").Value = 73.41986342149652
$ws.Range("G24").Value = 3.833625958456372
$ws.Range("J24").Value = 11.56254093473491
$ws.Range("K24").Value = 16.71532982942677
$ws.Range("L24").Value = 10.92584221621257
$ws.Range("M24").Value = 18.47012605573893
$ws.Range("B25").Value = 20.55043554600771
$ws.Range("C25").Value = 3.642648163456919
$ws.Range("D25").Value = 4.257944247699928
$ws.Range("F25").Value = 71.80480133351483
$ws.Range("G25").Value = 3.839850952165561
$ws.Range("J25").Value = 11.54079128903177
$ws.Range("K25").Value = 16.74202317816166
$ws.Range("L25").Value = 10.99525889667018
$ws.Range("M25").Value = 18.5667334342666
